$wb = $excel.ActiveWorkbook

# --- Add the new "SC" worksheet at the end of the tab order ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws9.Name = "SC"

$ws8 = $wb.Worksheets.Item("CO DR-0100 XML")

# --- Populate "CO DR-0100 XML" worksheet content ---
$ws8.Range("A1").Value = "Inputs that we need"
$ws8.Range("B2").Value = "CO DR-0100 XML tool.xlsm"
$ws8.Range("F2").Value = "XLSM file with macros that will help us to get state balancing"
$ws8.Range("B3").Value = "Detail worksheet"
$ws8.Range("F3").Value = "This detail file is downloaded from Taxsolver"
$ws8.Range("A6").Value = "Step by step"
$ws8.Range("B8").Value = 1
$ws8.Range("C8").Value = "Go to the detail tab"
$ws8.Range("B9").Value = 2
$ws8.Range("C9").Value = "Go to DR 0100 XML form page"
$ws8.Range("B10").Value = 3
$ws8.Range("C10").Value = "click on `"Export Detail as Worksheet`""
$ws8.Range("C11").Value = 3.1
$ws8.Range("D11").Value = "Save the file in some folder in the P drive"
$ws8.Range("C12").Value = 3.2
$ws8.Range("D12").Value = "This will be just a temp file"
$ws8.Range("B13").Value = 4
$ws8.Range("C13").Value = "Copy and paste the Detail worksheet file into the Tool file (in the Detail worksheet)"
$ws8.Range("B15").Value = "2nd part"
$ws8.Range("B7").Value = "1st part"
$ws8.Range("B16").Value = 1
$ws8.Range("C16").Value = "Go to Information Sheet in Taxsolver"
$ws8.Range("B17").Value = 2
$ws8.Range("C17").Value = "Get the State Registration ID"
$ws8.Range("G17").Value = "State ID"
$ws8.Range("B18").Value = 3
$ws8.Range("C18").Value = "Open browser"
$ws8.Range("G18").Value = "https://www.colorado.gov/revenueonline/"
$ws8.Hyperlinks.Add($ws8.Range("G18"), "https://www.colorado.gov/revenueonline/", [System.Reflection.Missing]::Value, "https://www.colorado.gov/revenueonline/") | Out-Null
$ws8.Range("B19").Value = 4
$ws8.Range("C19").Value = "Click on Sales and Use Tax --- Find Sales and Use Tax rates"
$ws8.Range("I19").Value = "At the bottom of the page"
$ws8.Range("B20").Value = 5
$ws8.Range("C20").Value = "Click on View Business Location rates"
$ws8.Range("B21").Value = 6
$ws8.Range("C21").Value = "Type the State ID (2) into the Colorado Account Number field"
$ws8.Range("B22").Value = 7
$ws8.Range("C22").Value = "Wait for the page to load and then, click on Export"
$ws8.Range("B23").Value = 8
$ws8.Range("C23").Value = "Download the file"
$ws8.Range("C24").Value = 8.1
$ws8.Range("D24").Value = "Check which browser to use and what configuration each user needs to do before running the bot"
$ws8.Range("C25").Value = 8.2
$ws8.Range("D25").Value = "You can save the file in any location in your local drive"
$ws8.Range("C26").Value = 8.3
$ws8.Range("D26").Value = "The output is a TXT file, which contains the colorado table info"
$ws8.Range("B27").Value = 9
$ws8.Range("C27").Value = "Get the text from the TXT file and split by line and by semicolon (;)"
$ws8.Range("C28").Value = 9.1
$ws8.Range("D28").Value = "You must have at the end 19 columns"
$ws8.Range("B29").Value = 10
$ws8.Range("C29").Value = "Paste the result datatable into the Tool file, Website worksheet."
$ws8.Range("B31").Value = "3hd part - Use XLSM tool file"
$ws8.Range("B32").Value = 1
$ws8.Range("C32").Value = "Click on Unhide Columns button"
$ws8.Range("B33").Value = 2
$ws8.Range("C33").Value = "Click on Hide Columns button"
$ws8.Range("B34").Value = 3
$ws8.Range("C34").Value = "In the Rate-Services Fees worksheet, check if the sum of each column is equals to 0"
$ws8.Range("C35").Value = "3.1 We can focus on the row 60, which starts with `"Level County`""
$ws8.Range("C36").Value = "3.1 If some of them are different to 0, the bot should mark this as a failure"
$ws8.Range("B37").Value = 4
$ws8.Range("C37").Value = "In the Balancing Sheet, check if TAX Not balanced and EXCEPT Not balanced cells are equals to 0"
$ws8.Range("C38").Value = "4.1 If some of them are different to 0, the bot should mark this as a failure"
$ws8.Range("M35").Value = "Ask Jay"

# --- Populate "SC" worksheet content ---
$ws9.Range("A1").Value = "Inputs that we need"
$ws9.Range("B2").Value = "SC ST-389 tool (2).xlsx"
$ws9.Range("B3").Value = "Detail worksheet"
$ws9.Range("F3").Value = "This detail file is downloaded from Taxsolver"
$ws9.Range("A5").Value = "Step by step"
$ws9.Range("B6").Value = "1st part"
$ws9.Range("B7").Value = 1
$ws9.Range("C7").Value = "Go to the detail tab"
$ws9.Range("B8").Value = 2
$ws9.Range("C8").Value = "Go to SC ST 389 form page"
$ws9.Range("B9").Value = 3
$ws9.Range("C9").Value = "click on `"Export Detail as Worksheet`""
$ws9.Range("C10").Value = 3.1
$ws9.Range("D10").Value = "Save the file in some folder in the P drive"
$ws9.Range("C11").Value = 3.2
$ws9.Range("D11").Value = "This will be just a temp file"
$ws9.Range("B12").Value = 4
$ws9.Range("C12").Value = "Copy and paste the Detail worksheet file into the Tool file (in the Detail worksheet)"
$ws9.Range("F2").Value = "XLSX file with macros and formulas that will help us to get state balancing"

# --- Selections (also drives which sheet ends up tabSelected) ---
$wsLA = $wb.Worksheets.Item("LA")
$wsLA.Range("B3").Select() | Out-Null

$wsCA = $wb.Worksheets.Item("CA")
$wsCA.Range("R26").Select() | Out-Null

$ws9.Range("H7").Select() | Out-Null

$ws8.Range("F25").Select() | Out-Null

